# Update LR-pairs data for Efna5-Epha2 with refreshed TPM values.
# New data drops the "ECs" sending-cluster rows and recomputes the
# remaining FAPs/MuSCs rows against the updated TPM inputs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Overwrite data rows 2-7 with the refreshed values ---
$ws.Range('A2').Value2 = 'FAPs'
$ws.Range('B2').Value2 = 'Efna5'
$ws.Range('C2').Value2 = 'Epha2'
$ws.Range('D2').Value2 = 'ECs'
$ws.Range('E2').Value2 = 3
$ws.Range('F2').Value2 = 1
$ws.Range('G2').Value2 = 2.900731333333333
$ws.Range('H2').Value2 = 8.702194
$ws.Range('I2').Value2 = 0.8130494232775288
$ws.Range('J2').Value2 = 0.8130494232775289
$ws.Range('K2').Value2 = 3
$ws.Range('L2').Value2 = 1
$ws.Range('M2').Value2 = 7.757543333333333
$ws.Range('N2').Value2 = 23.27263
$ws.Range('O2').Value2 = 0.4040769763164727
$ws.Range('P2').Value2 = 0.4040769763164727
$ws.Range('Q2').Value2 = 22.50254901669111
$ws.Range('R2').Value2 = 202.52294115022
$ws.Range('S2').Value2 = 0.3285345525538358
$ws.Range('T2').Value2 = 0.3285345525538358
$ws.Range('A3').Value2 = 'FAPs'
$ws.Range('B3').Value2 = 'Efna5'
$ws.Range('C3').Value2 = 'Epha2'
$ws.Range('D3').Value2 = 'FAPs'
$ws.Range('E3').Value2 = 3
$ws.Range('F3').Value2 = 1
$ws.Range('G3').Value2 = 2.900731333333333
$ws.Range('H3').Value2 = 8.702194
$ws.Range('I3').Value2 = 0.8130494232775288
$ws.Range('J3').Value2 = 0.8130494232775289
$ws.Range('K3').Value2 = 2
$ws.Range('L3').Value2 = 0.6666666666666666
$ws.Range('M3').Value2 = 0.312365
$ws.Range('N3').Value2 = 0.937095
$ws.Range('O3').Value2 = 0.01627055103446774
$ws.Range('P3').Value2 = 0.01627055103446774
$ws.Range('Q3').Value2 = 0.9060869429366667
$ws.Range('R3').Value2 = 8.154782486430001
$ws.Range('S3').Value2 = 0.0132287621349816
$ws.Range('T3').Value2 = 0.0132287621349816
$ws.Range('A4').Value2 = 'FAPs'
$ws.Range('B4').Value2 = 'Efna5'
$ws.Range('C4').Value2 = 'Epha2'
$ws.Range('D4').Value2 = 'MuSCs'
$ws.Range('E4').Value2 = 3
$ws.Range('F4').Value2 = 1
$ws.Range('G4').Value2 = 2.900731333333333
$ws.Range('H4').Value2 = 8.702194
$ws.Range('I4').Value2 = 0.8130494232775288
$ws.Range('J4').Value2 = 0.8130494232775289
$ws.Range('K4').Value2 = 3
$ws.Range('L4').Value2 = 1
$ws.Range('M4').Value2 = 11.12827366666667
$ws.Range('N4').Value2 = 33.384821
$ws.Range('O4').Value2 = 0.5796524726490594
$ws.Range('P4').Value2 = 0.5796524726490595
$ws.Range('Q4').Value2 = 32.28013211080822
$ws.Range('R4').Value2 = 290.521188997274
$ws.Range('S4').Value2 = 0.4712861085887113
$ws.Range('T4').Value2 = 0.4712861085887114
$ws.Range('A5').Value2 = 'MuSCs'
$ws.Range('B5').Value2 = 'Efna5'
$ws.Range('C5').Value2 = 'Epha2'
$ws.Range('D5').Value2 = 'ECs'
$ws.Range('E5').Value2 = 3
$ws.Range('F5').Value2 = 1
$ws.Range('G5').Value2 = 0.6669870000000001
$ws.Range('H5').Value2 = 2.000961
$ws.Range('I5').Value2 = 0.1869505767224711
$ws.Range('J5').Value2 = 0.1869505767224711
$ws.Range('K5').Value2 = 3
$ws.Range('L5').Value2 = 1
$ws.Range('M5').Value2 = 7.757543333333333
$ws.Range('N5').Value2 = 23.27263
$ws.Range('O5').Value2 = 0.4040769763164727
$ws.Range('P5').Value2 = 0.4040769763164727
$ws.Range('Q5').Value2 = 5.174180555270001
$ws.Range('R5').Value2 = 46.56762499743
$ws.Range('S5').Value2 = 0.07554242376263687
$ws.Range('T5').Value2 = 0.07554242376263687
$ws.Range('A6').Value2 = 'MuSCs'
$ws.Range('B6').Value2 = 'Efna5'
$ws.Range('C6').Value2 = 'Epha2'
$ws.Range('D6').Value2 = 'FAPs'
$ws.Range('E6').Value2 = 3
$ws.Range('F6').Value2 = 1
$ws.Range('G6').Value2 = 0.6669870000000001
$ws.Range('H6').Value2 = 2.000961
$ws.Range('I6').Value2 = 0.1869505767224711
$ws.Range('J6').Value2 = 0.1869505767224711
$ws.Range('K6').Value2 = 2
$ws.Range('L6').Value2 = 0.6666666666666666
$ws.Range('M6').Value2 = 0.312365
$ws.Range('N6').Value2 = 0.937095
$ws.Range('O6').Value2 = 0.01627055103446774
$ws.Range('P6').Value2 = 0.01627055103446774
$ws.Range('Q6').Value2 = 0.208343394255
$ws.Range('R6').Value2 = 1.875090548295
$ws.Range('S6').Value2 = 0.003041788899486143
$ws.Range('T6').Value2 = 0.003041788899486143
$ws.Range('A7').Value2 = 'MuSCs'
$ws.Range('B7').Value2 = 'Efna5'
$ws.Range('C7').Value2 = 'Epha2'
$ws.Range('D7').Value2 = 'MuSCs'
$ws.Range('E7').Value2 = 3
$ws.Range('F7').Value2 = 1
$ws.Range('G7').Value2 = 0.6669870000000001
$ws.Range('H7').Value2 = 2.000961
$ws.Range('I7').Value2 = 0.1869505767224711
$ws.Range('J7').Value2 = 0.1869505767224711
$ws.Range('K7').Value2 = 3
$ws.Range('L7').Value2 = 1
$ws.Range('M7').Value2 = 11.12827366666667
$ws.Range('N7').Value2 = 33.384821
$ws.Range('O7').Value2 = 0.5796524726490594
$ws.Range('P7').Value2 = 0.5796524726490595
$ws.Range('Q7').Value2 = 7.422413868109001
$ws.Range('R7').Value2 = 66.80172481298101
$ws.Range('S7').Value2 = 0.1083663640603481
$ws.Range('T7').Value2 = 0.1083663640603481

# --- Remove the old trailing rows (former "ECs" sender rows 8-10) ---
$ws.Rows.Item(8).Delete() | Out-Null
$ws.Rows.Item(8).Delete() | Out-Null
$ws.Rows.Item(8).Delete() | Out-Null

